$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.041.68'
$ws.Range("E2").Value = '  +0.63%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.798.04'
$ws.Range("E3").Value = '  +1.91%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '429.03'
$ws.Range("E5").Value = '  +6.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.21'
$ws.Range("E6").Value = '  +9.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.625'
$ws.Range("E7").Value = '  +3.80%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.732'
$ws.Range("E9").Value = '  +2.34%  '

$ws.Range("E10").Value = '  -7.95%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000307'
$ws.Range("E11").Value = '  -14.82%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.49'
$ws.Range("E12").Value = '  +5.48%  '

$ws.Range("E13").Value = '  +9.05%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.402.14'
$ws.Range("E14").Value = '  +2.82%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.12'
$ws.Range("E15").Value = '  +5.14%  '

$ws.Range("E16").Value = '  +0.56%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.796.38'
$ws.Range("E17").Value = '  +1.96%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.91'
$ws.Range("E18").Value = '  +2.58%  '

$ws.Range("E19").Value = '  +6.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '66.173.67'
$ws.Range("E20").Value = '  +0.60%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '405.03'
$ws.Range("E21").Value = '  +0.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.11'
$ws.Range("E22").Value = '  +4.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.20'
$ws.Range("E23").Value = '  +6.54%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.70'
$ws.Range("E24").Value = '  -0.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '36.83'
$ws.Range("E25").Value = '  +2.18%  '

$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.76'
$ws.Range("E26").Value = '  +38.70%  '

$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.27'
$ws.Range("E27").Value = '  +6.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.88'
$ws.Range("E28").Value = '  +8.12%  '

$ws.Range("E29").Value = '  -1.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.75'
$ws.Range("E30").Value = '  +11.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '707.37'
$ws.Range("E31").Value = '  +3.09%  '

$ws.Range("E32").Value = '  +13.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.72'
$ws.Range("E33").Value = '  +0.48%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '41.18'
$ws.Range("E34").Value = '  +7.50%  '

$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.61'
$ws.Range("E36").Value = '  +33.03%  '

$ws.Range("E37").Value = '  -3.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.58'
$ws.Range("E38").Value = '  +2.90%  '

$ws.Range("E39").Value = '  +5.14%  '

$ws.Range("E40").Value = '  +41.27%  '

$ws.Range("E41").Value = '  +2.18%  '

$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.140'
$ws.Range("E42").Value = '  +4.85%  '

$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0674'
$ws.Range("E43").Value = '  -8.43%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.37'
$ws.Range("E45").Value = '  +5.09%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.321'
$ws.Range("E46").Value = '  +12.02%  '

$ws.Range("E47").Value = '  -0.79%  '

$ws.Range("E48").Value = '  +4.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.06'
$ws.Range("E49").Value = '  +2.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.90'
$ws.Range("E50").Value = '  -2.00%  '

$ws.Range("E51").Value = '  +0.57%  '
